# Update cryptos worksheet values per Tue May  9 14:20:13 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.614.35"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.22"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.09"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4264"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.84"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07325"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8888"
$ws.Range("E11").Value = "  -4.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.85"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.858.41"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.572"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.354"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06931"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.03"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008882"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.44"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.628.34"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.994"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.69"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.087.46"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.983"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.61"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.93"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.65"
$ws.Range("E29").Value = "  +7.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.228"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.906"
$ws.Range("E31").Value = "  +11.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08925"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7638"
$ws.Range("E33").Value = "  -6.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.586"
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.969"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.102"
$ws.Range("E36").Value = "  -6.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05386"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.094"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01949"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.802"
$ws.Range("E41").Value = "  -5.98%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.911"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5108"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1657"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.276"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06590"
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4760"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.42"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.45"
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.631"
$ws.Range("E51").Value = "  -2.60%  "
